{"js": "// Apply the diff: update the date line and the 25 division-expression table cells.\n// Every old text string is unique within the document body, so a simple\n// search-and-replace (one pair at a time, oldest text -> newest text) is safe.\nconst replacements = [\n  [\"2025-10-10 Friday\", \"2025-10-11 Saturday\"],\n  [\"68\u00f76=11, 2\", \"22\u00f76=3, 4\"],\n  [\"95\u00f75=19, 0\", \"29\u00f76=4, 5\"],\n  [\"42\u00f77=6, 0\", \"48\u00f77=6, 6\"],\n  [\"75\u00f72=37, 1\", \"21\u00f72=10, 1\"],\n  [\"83\u00f78=10, 3\", \"69\u00f72=34, 1\"],\n  [\"35\u00f74=8, 3\", \"55\u00f79=6, 1\"],\n  [\"41\u00f77=5, 6\", \"66\u00f73=22, 0\"],\n  [\"60\u00f75=12, 0\", \"52\u00f73=17, 1\"],\n  [\"31\u00f78=3, 7\", \"21\u00f72=10, 1\"],\n  [\"21\u00f74=5, 1\", \"67\u00f73=22, 1\"],\n  [\"16\u00f73=5, 1\", \"28\u00f74=7, 0\"],\n  [\"29\u00f75=5, 4\", \"87\u00f76=14, 3\"],\n  [\"98\u00f77=14, 0\", \"51\u00f76=8, 3\"],\n  [\"72\u00f79=8, 0\", \"40\u00f73=13, 1\"],\n  [\"84\u00f77=12, 0\", \"68\u00f72=34, 0\"],\n  [\"62\u00f78=7, 6\", \"99\u00f77=14, 1\"],\n  [\"74\u00f75=14, 4\", \"36\u00f78=4, 4\"],\n  [\"60\u00f77=8, 4\", \"83\u00f75=16, 3\"],\n  [\"26\u00f74=6, 2\", \"48\u00f77=6, 6\"],\n  [\"93\u00f77=13, 2\", \"68\u00f73=22, 2\"],\n  [\"67\u00f72=33, 1\", \"91\u00f73=30, 1\"],\n  [\"90\u00f76=15, 0\", \"30\u00f72=15, 0\"],\n  [\"51\u00f79=5, 6\", \"78\u00f77=11, 1\"],\n  [\"78\u00f79=8, 6\", \"99\u00f72=49, 1\"],\n  [\"15\u00f75=3, 0\", \"16\u00f78=2, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the diff: update the date line and 25 division-expression cells.\n# Each old text is unique in the document, so Find/Replace targets exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2025-10-10 Friday\", \"2025-10-11 Saturday\")\n  ,@(\"68\u00f76=11, 2\", \"22\u00f76=3, 4\")\n  ,@(\"95\u00f75=19, 0\", \"29\u00f76=4, 5\")\n  ,@(\"42\u00f77=6, 0\", \"48\u00f77=6, 6\")\n  ,@(\"75\u00f72=37, 1\", \"21\u00f72=10, 1\")\n  ,@(\"83\u00f78=10, 3\", \"69\u00f72=34, 1\")\n  ,@(\"35\u00f74=8, 3\", \"55\u00f79=6, 1\")\n  ,@(\"41\u00f77=5, 6\", \"66\u00f73=22, 0\")\n  ,@(\"60\u00f75=12, 0\", \"52\u00f73=17, 1\")\n  ,@(\"31\u00f78=3, 7\", \"21\u00f72=10, 1\")\n  ,@(\"21\u00f74=5, 1\", \"67\u00f73=22, 1\")\n  ,@(\"16\u00f73=5, 1\", \"28\u00f74=7, 0\")\n  ,@(\"29\u00f75=5, 4\", \"87\u00f76=14, 3\")\n  ,@(\"98\u00f77=14, 0\", \"51\u00f76=8, 3\")\n  ,@(\"72\u00f79=8, 0\", \"40\u00f73=13, 1\")\n  ,@(\"84\u00f77=12, 0\", \"68\u00f72=34, 0\")\n  ,@(\"62\u00f78=7, 6\", \"99\u00f77=14, 1\")\n  ,@(\"74\u00f75=14, 4\", \"36\u00f78=4, 4\")\n  ,@(\"60\u00f77=8, 4\", \"83\u00f75=16, 3\")\n  ,@(\"26\u00f74=6, 2\", \"48\u00f77=6, 6\")\n  ,@(\"93\u00f77=13, 2\", \"68\u00f73=22, 2\")\n  ,@(\"67\u00f72=33, 1\", \"91\u00f73=30, 1\")\n  ,@(\"90\u00f76=15, 0\", \"30\u00f72=15, 0\")\n  ,@(\"51\u00f79=5, 6\", \"78\u00f77=11, 1\")\n  ,@(\"78\u00f79=8, 6\", \"99\u00f72=49, 1\")\n  ,@(\"15\u00f75=3, 0\", \"16\u00f78=2, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $old\"\n  }\n}\n\nWrite-Output \"Done\""}
